$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.834.54"
Set-TextValue $ws.Range("E2") "  +1.04%  "
Set-TextValue $ws.Range("D3") "1.756.69"
Set-TextValue $ws.Range("E3") "  +0.13%  "
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "327.84"
Set-TextValue $ws.Range("E5") "  +1.07%  "
Set-TextValue $ws.Range("E6") "  +0.03%  "
Set-TextValue $ws.Range("D7") "0.4562"
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("E8") "  -1.76%  "
Set-TextValue $ws.Range("D9") "41.97"
Set-TextValue $ws.Range("E9") "  +1.02%  "
Set-TextValue $ws.Range("D10") "0.07350"
Set-TextValue $ws.Range("E10") "  -1.61%  "
Set-TextValue $ws.Range("D11") "1.084"
Set-TextValue $ws.Range("E11") "  -0.37%  "
Set-TextValue $ws.Range("E12") "  +0.11%  "
Set-TextValue $ws.Range("D13") "20.61"
Set-TextValue $ws.Range("E13") "  -0.90%  "
Set-TextValue $ws.Range("D14") "5.978"
Set-TextValue $ws.Range("E14") "  -0.69%  "
Set-TextValue $ws.Range("D15") "7.170"
Set-TextValue $ws.Range("E15") "  -0.07%  "
Set-TextValue $ws.Range("D16") "1.757.08"
Set-TextValue $ws.Range("E16") "  -0.21%  "
Set-TextValue $ws.Range("D17") "91.57"
Set-TextValue $ws.Range("E17") "  -2.59%  "
Set-TextValue $ws.Range("E18") "  -0.19%  "
Set-TextValue $ws.Range("D19") "0.06412"
Set-TextValue $ws.Range("E19") "  +0.20%  "
Set-TextValue $ws.Range("E20") "  +0.06%  "
Set-TextValue $ws.Range("E21") "  -1.71%  "
Set-TextValue $ws.Range("D22") "5.745"
Set-TextValue $ws.Range("E22") "  +0.00%  "
Set-TextValue $ws.Range("D23") "27.865.72"
Set-TextValue $ws.Range("E23") "  +0.96%  "
Set-TextValue $ws.Range("E24") "  -0.45%  "
Set-TextValue $ws.Range("E25") "  +3.64%  "
Set-TextValue $ws.Range("D26") "162.47"
Set-TextValue $ws.Range("E26") "  -1.94%  "
Set-TextValue $ws.Range("D27") "20.02"
Set-TextValue $ws.Range("E27") "  -0.56%  "
Set-TextValue $ws.Range("D28") "1.959.72"
Set-TextValue $ws.Range("E28") "  +0.00%  "
Set-TextValue $ws.Range("D29") "2.159"
Set-TextValue $ws.Range("E29") "  +1.15%  "
Set-TextValue $ws.Range("D30") "123.32"
Set-TextValue $ws.Range("E30") "  -1.90%  "
Set-TextValue $ws.Range("D31") "1.080"
Set-TextValue $ws.Range("E31") "  -1.00%  "
Set-TextValue $ws.Range("D32") "0.09269"
Set-TextValue $ws.Range("E32") "  +0.45%  "
Set-TextValue $ws.Range("D33") "3.658"
Set-TextValue $ws.Range("E33") "  -0.03%  "
Set-TextValue $ws.Range("D34") "5.524"
Set-TextValue $ws.Range("E34") "  -0.19%  "
Set-TextValue $ws.Range("D35") "11.73"
Set-TextValue $ws.Range("E35") "  +0.00%  "
Set-TextValue $ws.Range("D36") "0.06092"
Set-TextValue $ws.Range("E36") "  +1.24%  "
Set-TextValue $ws.Range("D37") "0.02250"
Set-TextValue $ws.Range("E37") "  -1.45%  "
Set-TextValue $ws.Range("D38") "0.2060"
Set-TextValue $ws.Range("E38") "  -1.63%  "
Set-TextValue $ws.Range("D39") "4.894"
Set-TextValue $ws.Range("E39") "  -0.58%  "
Set-TextValue $ws.Range("D40") "0.6188"
Set-TextValue $ws.Range("E40") "  -1.83%  "
Set-TextValue $ws.Range("D41") "1.176"
Set-TextValue $ws.Range("E41") "  -0.41%  "
Set-TextValue $ws.Range("D42") "1.368"
Set-TextValue $ws.Range("E42") "  -1.49%  "
Set-TextValue $ws.Range("D43") "7.752"
Set-TextValue $ws.Range("E43") "  -0.62%  "
Set-TextValue $ws.Range("D44") "13.02"
Set-TextValue $ws.Range("E44") "  -0.86%  "
Set-TextValue $ws.Range("D45") "3.728"
Set-TextValue $ws.Range("E45") "  +0.32%  "
Set-TextValue $ws.Range("D46") "0.5802"
Set-TextValue $ws.Range("E46") "  -1.03%  "
Set-TextValue $ws.Range("D47") "122.30"
Set-TextValue $ws.Range("D48") "1.925"
Set-TextValue $ws.Range("E48") "  -0.38%  "
Set-TextValue $ws.Range("D49") "0.06784"
Set-TextValue $ws.Range("E49") "  -1.56%  "
Set-TextValue $ws.Range("D50") "1.118"
Set-TextValue $ws.Range("E50") "  -1.12%  "
Set-TextValue $ws.Range("D51") "72.21"
Set-TextValue $ws.Range("E51") "  +0.07%  "
